$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetText($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

SetText "D2" '67.719.56'
SetText "E2" '  -1.14%  '
SetText "D3" '3.795.69'
SetText "E3" '  +1.47%  '
SetText "E4" '  -0.11%  '
SetText "D5" '595.48'
SetText "E5" '  +0.12%  '
SetText "D6" '166.72'
SetText "E6" '  -0.41%  '
SetText "D7" '3.795.50'
SetText "E7" '  +1.53%  '
SetText "E8" '  +0.09%  '
SetText "D9" '0.519'
SetText "E9" '  -0.14%  '
SetText "E10" '  +0.00%  '
SetText "D11" '6.35'
SetText "E11" '  -1.93%  '
SetText "D12" '0.451'
SetText "E12" '  +0.48%  '
SetText "E13" '  -0.80%  '
SetText "D14" '36.37'
SetText "E14" '  +0.24%  '
SetText "D15" '4.433.45'
SetText "E15" '  +1.43%  '
SetText "D16" '3.796.74'
SetText "E16" '  +1.74%  '
SetText "D17" '18.69'
SetText "E17" '  +4.38%  '
SetText "D18" '67.678.31'
SetText "E18" '  -1.15%  '
SetText "E19" '  +0.36%  '
SetText "D20" '7.00'
SetText "E20" '  -0.32%  '
SetText "D21" '10.25'
SetText "E21" '  -3.98%  '
SetText "D22" '458.59'
SetText "E22" '  -1.90%  '
SetText "E23" '  -0.28%  '
SetText "D24" '0.0000156'
SetText "E24" '  +8.69%  '
SetText "D25" '83.61'
SetText "E25" '  -0.72%  '
SetText "D26" '11.91'
SetText "E26" '  -1.35%  '
SetText "E27" '  -2.71%  '
SetText "E28" '  -0.49%  '
SetText "E29" '  +0.14%  '
SetText "E30" '  +0.04%  '
SetText "E31" '  -0.49%  '
SetText "B32" 'ImmutableX'
SetText "C32" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
SetText "D32" '2.20'
SetText "E32" '  +0.70%  '
SetText "B33" 'EthereumClassic'
SetText "C33" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
SetText "D33" '29.86'
SetText "E33" '  -0.29%  '
SetText "E34" '  -0.66%  '
SetText "E35" '  -0.06%  '
SetText "D36" '3.749.05'
SetText "E36" '  +1.38%  '
SetText "E37" '  -1.40%  '
SetText "D38" '3.36'
SetText "E38" '  -2.65%  '
SetText "E39" '  -0.38%  '
SetText "E40" '  +0.45%  '
SetText "E41" '  -0.69%  '
SetText "E42" '  -0.03%  '
SetText "D44" '44.72'
SetText "E44" '  +3.98%  '
SetText "E45" '  -2.10%  '
SetText "D46" '47.11'
SetText "E46" '  +2.61%  '
SetText "E47" '  -2.58%  '
SetText "D48" '148.35'
SetText "E48" '  +0.91%  '
SetText "D49" '394.00'
SetText "E49" '  +0.71%  '
SetText "E50" '  -4.99%  '
SetText "D51" '2.760.90'
SetText "E51" '  +2.41%  '
